# Insert a new data row at row 255 (pushing existing rows 255-352 down to 256-353)
# and populate it with the new weekly price observation, per the commit
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 255..352 down by one, creating a blank row 255.
$ws.Rows(255).Insert()

# Populate the newly inserted row 255 with the new record.
$ws.Range("A255").Value = 7
$ws.Range("B255").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C255").Value = "Ñuble"
$ws.Range("D255").Value = 44900
$ws.Range("E255").Value = 16
$ws.Range("F255").Value = 100114013
$ws.Range("G255").Value = "Zanahoria"
$ws.Range("H255").Value = "Sin especificar"
$ws.Range("I255").Value = "Primera"
$ws.Range("J255").Value = 160
$ws.Range("K255").Value = 8500
$ws.Range("L255").Value = 9000
$ws.Range("M255").Value = 8750
$ws.Range("N255").Value = "$/saco 20 kilos"
$ws.Range("O255").Value = "Región de Ñuble"
$ws.Range("P255").Value = 438
$ws.Range("Q255").Value = 20
$ws.Range("R255").Value = "Hortaliza"
